$d = $word.ActiveDocument

function Get-ParaText($p) {
    $t = $p.Range.Text
    if ($t.Length -gt 0) {
        return $t.Substring(0, $t.Length - 1)
    }
    return $t
}

# --- 1. Merge the three runs that make up
#        "Variables may be set to " + "an expression" + "(such as "
#     into a single run by doing a no-op Find/Replace over that exact span.
#     (Replacing a range's text with itself normalises it into one run.)
$d.Content.Find.Execute("an expression(such as ", $false, $false, $false, $false, $false, $true, 1, $false, "an expression(such as ", 2)

# --- 2. Remove the old bookmark that wrapped "c = a + b" ...
$d.Bookmarks("_GoBack").Delete()

# --- 3. Add the missing w:pos="840" tab stop to the five sub-bullet
#        paragraphs that are still only using the w:pos="425" stop.
$targets = @(
    "Input, Read, Fetch, Get > single variable",
    "New variables may be found ",
    "Output, Print, Show, Display > single variable (or String, with double quote)",
    "Two variables (or one variable and one value) separated by one relational operator",
    "At least two simple conditions separated by logical operators (and possibly by parentheses)"
)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = Get-ParaText $p
    foreach ($target in $targets) {
        if ($t.Equals($target)) {
            $p.Format.TabStops.Add(42)
        }
    }
}

# --- 4. Re-create "_GoBack" as an empty bookmark at the very end of the
#        "Output, Print, Show, Display ..." paragraph's text (just before
#        its paragraph mark).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = Get-ParaText $p
    if ($t.Equals("Output, Print, Show, Display > single variable (or String, with double quote)")) {
        $pos = $p.Range.End - 1
        $r = $d.Range($pos, $pos)
        $r.InsertAfter("X")
        $r2 = $d.Range($pos, $pos + 1)
        $d.Bookmarks.Add("_GoBack", $r2)
        $r3 = $d.Range($pos, $pos + 1)
        $r3.Text = ""
    }
}

# --- 5. Delete the now-duplicate "New variables may be found" paragraph
#        that immediately follows the "Output, Print, Show, Display ..."
#        paragraph (the one with no trailing space).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = Get-ParaText $p
    if ($t.Equals("New variables may be found")) {
        $p.Range.Delete()
        break
    }
}
